# Adds a new weekly price-observation row for "Fruta, Terminal La Palmera de
# La Serena - Ciruela" at the top of the existing data block (row 312),
# pushing the previously-existing rows 312..373 down by one row (to
# 313..374). This mirrors the "Fruta / hortaliza, semanal" weekly-refresh
# commit: a brand-new record is inserted and everything below shifts down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 312; Excel shifts 312..373 -> 313..374 automatically
# and the worksheet dimension grows from A1:T373 to A1:T374 on save.
$ws.Rows(312).Insert()

# Populate the newly-inserted row 312 with the new observation.
$ws.Cells.Item(312, 1).Value2 = 8
$ws.Cells.Item(312, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(312, 3).Value = "Coquimbo"
$ws.Cells.Item(312, 4).Value2 = 45005
$ws.Cells.Item(312, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(312, 5).Value2 = 4
$ws.Cells.Item(312, 6).Value = "Fruta"
$ws.Cells.Item(312, 7).Value2 = 100103
$ws.Cells.Item(312, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(312, 9).Value2 = 100103002
$ws.Cells.Item(312, 10).Value = "Ciruela"
$ws.Cells.Item(312, 11).Value = "Angeleno"
$ws.Cells.Item(312, 12).Value = "Primera"
$ws.Cells.Item(312, 13).Value2 = 18
$ws.Cells.Item(312, 14).Value2 = 200000
$ws.Cells.Item(312, 15).Value2 = 210000
$ws.Cells.Item(312, 16).Value2 = 205000
$ws.Cells.Item(312, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(312, 18).Value = "Región Metropolitana"
$ws.Cells.Item(312, 19).Value2 = 456
$ws.Cells.Item(312, 20).Value2 = 450
